$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.509722597846064
$ws.Range("D2").Value = 0.01342731672374953
$ws.Range("E2").Value = 0.6436769344872033
$ws.Range("F2").Value = 0.3136141733259663
$ws.Range("G2").Value = 0.1832981381767809
$ws.Range("H2").Value = 0.3176921567105495
$ws.Range("I2").Value = 0.8116579071855838
$ws.Range("L2").Value = 0.4928928701960729
$ws.Range("O2").Value = 0.9059872803690894
$ws.Range("B3").Value = 1.367318995150583
$ws.Range("D3").Value = 0.01172222869912076
$ws.Range("E3").Value = 0.611446073590713
$ws.Range("F3").Value = 0.304678709452638
$ws.Range("G3").Value = 0.1751760899433208
$ws.Range("H3").Value = 0.3178574265688781
$ws.Range("I3").Value = 0.7992033266433225
$ws.Range("L3").Value = 0.4354744488040296
$ws.Range("O3").Value = 0.8888264814999332
$ws.Range("B4").Value = 1.279574098457772
$ws.Range("D4").Value = 0.01067205223134948
$ws.Range("E4").Value = 0.5918752725521017
$ws.Range("F4").Value = 0.2994982668508115
$ws.Range("G4").Value = 0.1704113415099826
$ws.Range("H4").Value = 0.3181992942576599
$ws.Range("I4").Value = 0.7922808458555934
$ws.Range("L4").Value = 0.4000915507226068
$ws.Range("O4").Value = 0.8792536921792191
$ws.Range("B5").Value = 1.24374216191012
$ws.Range("D5").Value = 0.01024331856779526
$ws.Range("E5").Value = 0.5839566661540232
$ws.Range("F5").Value = 0.2974638381134582
$ws.Range("G5").Value = 0.1685251413409148
$ws.Range("H5").Value = 0.3183989519331334
$ws.Range("I5").Value = 0.7896425648281138
$ws.Range("L5").Value = 0.3856416516871946
$ws.Range("O5").Value = 0.8755940245228118
$ws.Range("B6").Value = 1.23778782697309
$ws.Range("D6").Value = 0.01017208178477347
$ws.Range("E6").Value = 0.5826452556292026
$ws.Range("F6").Value = 0.297130643200461
$ws.Range("G6").Value = 0.1682152777668762
$ws.Range("H6").Value = 0.3184357464839138
$ws.Range("I6").Value = 0.7892155279701498
$ws.Range("L6").Value = 0.3832404107822924
$ws.Range("O6").Value = 0.8750008827939126
$ws.Range("B7").Value = 1.279091157620428
$ws.Range("D7").Value = 0.01066627327958258
$ws.Range("E7").Value = 0.5917682479839783
$ws.Range("F7").Value = 0.2994705198983283
$ws.Range("G7").Value = 0.170385679484518
$ws.Range("H7").Value = 0.3182017427209374
$ws.Range("I7").Value = 0.7922445247642358
$ws.Range("L7").Value = 0.3998967987793662
$ws.Range("O7").Value = 0.8792033609650503
$ws.Range("B8").Value = 1.460687249273406
$ws.Range("D8").Value = 0.01284009814737175
$ws.Range("E8").Value = 0.6325191865076647
$ws.Range("F8").Value = 0.3104695415365342
$ws.Range("G8").Value = 0.1804512874831516
$ws.Range("H8").Value = 0.3176991848945079
$ws.Range("I8").Value = 0.8072134212547724
$ws.Range("L8").Value = 0.4731221187219603
$ws.Range("O8").Value = 0.8998694969822623
$ws.Range("B9").Value = 1.814262127218115
$ws.Range("D9").Value = 0.01707574886289365
$ws.Range("E9").Value = 0.7141061464823224
$ws.Range("F9").Value = 0.3344810799574987
$ws.Range("G9").Value = 0.2019725385915621
$ws.Range("H9").Value = 0.3186262248095915
$ws.Range("I9").Value = 0.842301567055415
$ws.Range("L9").Value = 0.6156645770115574
$ws.Range("O9").Value = 0.9480969420248186
$ws.Range("B10").Value = 2.072399617091833
$ws.Range("D10").Value = 0.02016942592716475
$ws.Range("E10").Value = 0.7749904287297511
$ws.Range("F10").Value = 0.3536339400702531
$ws.Range("G10").Value = 0.218899324773318
$ws.Range("H10").Value = 0.3204810404215124
$ws.Range("I10").Value = 0.8715616011315319
$ws.Range("L10").Value = 0.7197082195258702
$ws.Range("O10").Value = 0.9883014117132234
$ws.Range("B11").Value = 2.189461777525821
$ws.Range("D11").Value = 0.02157252266362519
$ws.Range("E11").Value = 0.8028769334174797
$ws.Range("F11").Value = 0.362680228357064
$ws.Range("G11").Value = 0.2268480350558235
$ws.Range("H11").Value = 0.3215814053081374
$ws.Range("I11").Value = 0.8856261515798849
$ws.Range("L11").Value = 0.7668844441146518
$ws.Range("O11").Value = 1.007644274317755
$ws.Range("B12").Value = 2.233735674880677
$ws.Range("D12").Value = 0.02210319599090838
$ws.Range("E12").Value = 0.8134626968632688
$ws.Range("F12").Value = 0.3661541304324203
$ws.Range("G12").Value = 0.2298942092632501
$ws.Range("H12").Value = 0.3220351068897571
$ws.Range("I12").Value = 0.8910601847422157
$ws.Range("L12").Value = 0.784725888102713
$ws.Range("O12").Value = 1.015121649130037
$ws.Range("B13").Value = 2.224202982723341
$ws.Range("D13").Value = 0.02198893549758907
$ws.Range("E13").Value = 0.8111817428019492
$ws.Range("F13").Value = 0.3654038118802134
$ws.Range("G13").Value = 0.2292365458584982
$ws.Range("H13").Value = 0.3219357459063161
$ws.Range("I13").Value = 0.8898850660899171
$ws.Range("L13").Value = 0.7808844610343613
$ws.Range("O13").Value = 1.013504455373862
$ws.Range("B14").Value = 2.19310533355906
$ws.Range("D14").Value = 0.02161619471387866
$ws.Range("E14").Value = 0.803747323385906
$ws.Range("F14").Value = 0.3629650592586131
$ws.Range("G14").Value = 0.2270979183590924
$ws.Range("H14").Value = 0.3216179890539479
$ws.Range("I14").Value = 0.8860710482469898
$ws.Range("L14").Value = 0.7683527407895099
$ws.Range("O14").Value = 1.008256376567687
$ws.Range("B15").Value = 2.174049892465121
$ws.Range("D15").Value = 0.02138779451530581
$ws.Range("E15").Value = 0.7991968340471232
$ws.Range("F15").Value = 0.3614775478797014
$ws.Range("G15").Value = 0.2257926700122397
$ws.Range("H15").Value = 0.3214281780484072
$ws.Range("I15").Value = 0.8837489195686743
$ws.Range("L15").Value = 0.7606736519171307
$ws.Range("O15").Value = 1.005061690280058
$ws.Range("B16").Value = 2.064741761908806
$ws.Range("D16").Value = 0.02007764163531078
$ws.Range("E16").Value = 0.773171668052413
$ws.Range("F16").Value = 0.3530494819275418
$ws.Range("G16").Value = 0.2183848966997886
$ws.Range("H16").Value = 0.3204143025908905
$ws.Range("I16").Value = 0.8706575934293284
$ws.Range("L16").Value = 0.7166219575352386
$ws.Range("O16").Value = 0.9870586098447802
$ws.Range("B17").Value = 1.997589339579974
$ws.Range("D17").Value = 0.01927279248182856
$ws.Range("E17").Value = 0.7572536014188813
$ws.Range("F17").Value = 0.3479647856020023
$ws.Range("G17").Value = 0.2139044459822657
$ws.Range("H17").Value = 0.3198581289088338
$ws.Range("I17").Value = 0.8628193664483632
$ws.Range("L17").Value = 0.6895575319411194
$ws.Range("O17").Value = 0.9762849812286163
$ws.Range("B18").Value = 1.958930744067743
$ws.Range("D18").Value = 0.0188094685535134
$ws.Range("E18").Value = 0.7481159208264927
$ws.Range("F18").Value = 0.3450715806814415
$ws.Range("G18").Value = 0.2113507867658768
$ws.Range("H18").Value = 0.3195623778985919
$ws.Range("I18").Value = 0.8583820155910047
$ws.Range("L18").Value = 0.6739763807733539
$ws.Range("O18").Value = 0.9701873818043225
$ws.Range("B19").Value = 1.945835790659714
$ws.Range("D19").Value = 0.01865252825518837
$ws.Range("E19").Value = 0.7450251935464109
$ws.Range("F19").Value = 0.344097371076181
$ws.Range("G19").Value = 0.2104901639726791
$ws.Range("H19").Value = 0.3194663850537154
$ws.Range("I19").Value = 0.8568918079231196
$ws.Range("L19").Value = 0.668698428419475
$ws.Range("O19").Value = 0.968139826231436
$ws.Range("B20").Value = 2.004741396769361
$ws.Range("D20").Value = 0.01935851132381572
$ws.Range("E20").Value = 0.7589462583068212
$ws.Range("F20").Value = 0.348502810495873
$ws.Range("G20").Value = 0.2143789750890335
$ws.Range("H20").Value = 0.319914834727669
$ws.Range("I20").Value = 0.863646414370848
$ws.Range("L20").Value = 0.6924400852925601
$ws.Range("O20").Value = 0.9774215869209684
$ws.Range("B21").Value = 2.202240978498821
$ws.Range("D21").Value = 0.02172569565068727
$ws.Range("E21").Value = 0.8059303061254894
$ws.Range("F21").Value = 0.363680067605614
$ws.Range("G21").Value = 0.2277251001405176
$ws.Range("H21").Value = 0.3217103163663779
$ws.Range("I21").Value = 0.88718838667792
$ws.Range("L21").Value = 0.7720342479461522
$ws.Range("O21").Value = 1.00979371254985
$ws.Range("B22").Value = 2.330996700224091
$ws.Range("D22").Value = 0.02326898626181162
$ws.Range("E22").Value = 0.8367864711651123
$ws.Range("F22").Value = 0.373880770071608
$ws.Range("G22").Value = 0.2366585828590502
$ws.Range("H22").Value = 0.3230995867041599
$ws.Range("I22").Value = 0.9032043358073878
$ws.Range("L22").Value = 0.8239181108587275
$ws.Range("O22").Value = 1.031841065799711
$ws.Range("B23").Value = 2.262307559162707
$ws.Range("D23").Value = 0.02244566410227122
$ws.Range("E23").Value = 0.820304799532849
$ws.Range("F23").Value = 0.3684106074625788
$ws.Range("G23").Value = 0.2318711697769089
$ws.Range("H23").Value = 0.3223383206137953
$ws.Range("I23").Value = 0.8945987904238706
$ws.Range("L23").Value = 0.7962394630709184
$ws.Range("O23").Value = 1.019992142171304
$ws.Range("B24").Value = 2.001508112795705
$ws.Range("D24").Value = 0.01931975971778854
$ws.Range("E24").Value = 0.7581809649315261
$ws.Range("F24").Value = 0.3482594758467243
$ws.Range("G24").Value = 0.2141643713251966
$ws.Range("H24").Value = 0.319889123274109
$ws.Range("I24").Value = 0.8632722912030744
$ws.Range("L24").Value = 0.6911369496687882
$ws.Range("O24").Value = 0.9769074275967569
$ws.Range("B25").Value = 1.718891457340476
$ws.Range("D25").Value = 0.01593297942372374
$ws.Range("E25").Value = 0.6918636263625899
$ws.Range("F25").Value = 0.3277214917820217
$ws.Range("G25").Value = 0.1959566812079885
$ws.Range("H25").Value = 0.3181698765240952
$ws.Range("I25").Value = 0.8321970357804958
$ws.Range("L25").Value = 0.5772199952485551
$ws.Range("O25").Value = 0.9342176897655747
